$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.604.06'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.099.29'
$ws.Range('E3').Value = '  +9.81%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''252.85'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('D6').Value = '''0.660'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '''49.58'
$ws.Range('E8').Value = '  +5.85%  '
$ws.Range('D9').Value = '''60.44'
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('D10').Value = '''0.377'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('D13').Value = '''14.68'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '2.406.05'
$ws.Range('E14').Value = '  +9.60%  '
$ws.Range('D15').Value = '''0.835'
$ws.Range('E15').Value = '  +2.11%  '
$ws.Range('D16').Value = '2.104.14'
$ws.Range('E16').Value = '  +10.04%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '36.567.72'
$ws.Range('D19').Value = '''73.08'
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('D20').Value = '0.0₃0836'
$ws.Range('E20').Value = '  -2.95%  '
$ws.Range('D21').Value = '''13.35'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('E22').Value = '  -4.49%  '
$ws.Range('D23').Value = '''5.26'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('D26').Value = '''171.09'
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').Value = '''21.22'
$ws.Range('E27').Value = '  +13.30%  '
$ws.Range('D28').Value = '''9.20'
$ws.Range('E28').Value = '  +4.15%  '
$ws.Range('E29').Value = '  -9.84%  '
$ws.Range('D30').Value = '''28.55'
$ws.Range('E30').Value = '  +53.03%  '
$ws.Range('E31').Value = '  -5.07%  '
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('D33').Value = '''0.0618'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').Value = '''2.43'
$ws.Range('E34').Value = '  +20.69%  '
$ws.Range('D35').Value = '''0.975'
$ws.Range('E35').Value = '  +11.46%  '
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('E39').Value = '  -5.70%  '
$ws.Range('E40').Value = '  -11.18%  '
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('E42').Value = '  +6.16%  '
$ws.Range('D43').Value = '''97.92'
$ws.Range('E43').Value = '  -7.36%  '
$ws.Range('D44').Value = '''16.48'
$ws.Range('E44').Value = '  -8.28%  '
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').Value = '1.339.79'
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('D47').Value = '''0.0846'
$ws.Range('E47').Value = '  +4.15%  '
$ws.Range('D48').Value = '''7.09'
$ws.Range('E48').Value = '  +9.42%  '
$ws.Range('D49').Value = '''2.89'
$ws.Range('E49').Value = '  +3.13%  '
$ws.Range('D50').Value = '2.299.20'
$ws.Range('E50').Value = '  +9.83%  '
$ws.Range('E51').Value = '  -6.39%  '
